$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    $r = $ws.Range($range)
    $r.NumberFormat = "@"
    $r.Value = $text
    $r.NumberFormat = "General"
    $r.Style = "Normal"
}

# Update the regression coefficient table (values are stored as text, even
# when they look numeric, e.g. "1.65", "0.1" - force text so Excel does not
# reinterpret them as numbers)
Set-TextValue "B2" "-0.29***"
Set-TextValue "C2" "-0.02***"
Set-TextValue "D2" "0.26*"
Set-TextValue "E2" "-3.7*"

Set-TextValue "B3" "-1.02*"
Set-TextValue "C3" "-0.42***"
Set-TextValue "D3" "1.65"
Set-TextValue "E3" "-1.49"

Set-TextValue "B4" "0.06***"
Set-TextValue "C4" "-0.0*"
Set-TextValue "D4" "0.32***"
Set-TextValue "E4" "0.53"

Set-TextValue "B5" "0.01***"
Set-TextValue "C5" "0.0***"
Set-TextValue "D5" "-0.01*"
Set-TextValue "E5" "0.1"
